$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (column D) cells retain their exact text representation
# (avoids Excel auto-converting numeric-looking strings into floating point numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.117.66"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.385.15"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "661.81"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.382.41"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.97"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.781.31"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.017.19"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.06"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.420.34"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.21"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.525"
$ws.Range("E21").Value = "  -6.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.95"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "511.85"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.00"
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000201"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.567.12"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("E33").Value = "  -5.93%  "
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.565"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.90"
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "535.07"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.63"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.62"
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  +5.78%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.97"
$ws.Range("E51").Value = "  +3.13%  "
